# Insert two new rows at row 7 (pushing existing rows 7.. down to 9..),
# then populate the two new rows with the data reported in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:A8").EntireRow.Insert()

# New row 7
$row7 = @(3, "Femacal de La Calera", "Coquimbo", 44473, 5, 100112026, "Haba", "Sin especificar", "Primera", 85, 9000, 9500, 9265, "`$/saco 25 kilos", "Provincia de Limarí", 371, 25, "Hortaliza")
for ($c = 0; $c -lt $row7.Length; $c++) {
    $ws.Cells.Item(7, $c + 1).Value = $row7[$c]
}

# New row 8
$row8 = @(3, "Femacal de La Calera", "Coquimbo", 44473, 5, 100112026, "Haba", "Sin especificar", "Segunda", 40, 8000, 8000, 8000, "`$/saco 25 kilos", "Provincia de Limarí", 320, 25, "Hortaliza")
for ($c = 0; $c -lt $row8.Length; $c++) {
    $ws.Cells.Item(8, $c + 1).Value = $row8[$c]
}
